# Insert a new data row at row 206 (pushing existing rows 206.. down by one)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("206:206").Insert()

$ws.Cells.Item(206, 1).Value2 = 4
$ws.Cells.Item(206, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(206, 3).Value2 = "Los Lagos"
$ws.Cells.Item(206, 4).Value2 = 44762
$ws.Cells.Item(206, 5).Value2 = 10
$ws.Cells.Item(206, 6).Value2 = "Fruta"
$ws.Cells.Item(206, 7).Value2 = 100102
$ws.Cells.Item(206, 8).Value2 = "Cítricos"
$ws.Cells.Item(206, 9).Value2 = 100102006
$ws.Cells.Item(206, 10).Value2 = "Pomelo"
$ws.Cells.Item(206, 11).Value2 = "Start Ruby"
$ws.Cells.Item(206, 12).Value2 = "Primera"
$ws.Cells.Item(206, 13).Value2 = 80
$ws.Cells.Item(206, 14).Value2 = 15000
$ws.Cells.Item(206, 15).Value2 = 15000
$ws.Cells.Item(206, 16).Value2 = 15000
$ws.Cells.Item(206, 17).Value2 = "$/caja 14 kilos empedrada"
$ws.Cells.Item(206, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(206, 19).Value2 = 1071
$ws.Cells.Item(206, 20).Value2 = 14
